$d = $word.ActiveDocument
$ellipsis = [char]0x2026

# 1. Replace "Next permutation" text with "Re-do a lot of array problems later…"
$rng = $d.Content
$found = $rng.Find.Execute("Next permutation", $true, $false, $false, $false, $false, $true, 1, $false, "Re-do a lot of array problems later$ellipsis", 2)

# Find the paragraph that now contains the replaced text, insert the new
# "Strings:" block of paragraphs right after it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Re-do a lot of array problems later$ellipsis") {
        $p.Range.InsertParagraphAfter()
        $p1 = $p.Next()
        $p1.Range.Text = "Strings:"

        $p1.Range.InsertParagraphAfter()
        $p2 = $p1.Next()
        $p2.Range.Text = "76 Minimum Window Substring"

        $p2.Range.InsertParagraphAfter()
        $p3 = $p2.Next()
        $p3.Range.Text = "49 Group Anagrams"

        $p3.Range.InsertParagraphAfter()
        $p4 = $p3.Next()
        $p4.Range.Text = "17 Letter Combinations in a Phone Number"
        break
    }
}

# 2. Insert "Graph:" and "Critical Edges" paragraphs right after "DP:" (the
# last paragraph of the body). Inserting at the very end of the document
# needs the content range collapsed-to-end + re-fetch-by-index pattern,
# since the trailing paragraph's "Next" reference doesn't reseat reliably.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$g1 = $d.Paragraphs($d.Paragraphs.Count)
$g1.Range.Text = "Graph:"

$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()
$g2 = $d.Paragraphs($d.Paragraphs.Count)
$g2.Range.Text = "Critical Edges"
